$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.692.78'
$ws.Range("E2").Value = '  -3.34%  '
$ws.Range("D3").Value = '3.309.39'
$ws.Range("E3").Value = '  -5.86%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''550.95'
$ws.Range("E5").Value = '  -4.51%  '
$ws.Range("D6").Value = '''172.53'
$ws.Range("E6").Value = '  -3.27%  '
$ws.Range("D7").Value = '''0.610'
$ws.Range("E7").Value = '  -4.04%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '3.303.32'
$ws.Range("E9").Value = '  -5.80%  '
$ws.Range("D10").Value = '''0.618'
$ws.Range("E10").Value = '  -2.42%  '
$ws.Range("D11").Value = '''0.159'
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("D12").Value = '''53.40'
$ws.Range("E12").Value = '  -3.08%  '
$ws.Range("D13").Value = '''0.0000270'
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").Value = '''9.00'
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").Value = '3.830.46'
$ws.Range("E15").Value = '  -6.08%  '
$ws.Range("D16").Value = '''18.18'
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("E17").Value = '  -3.56%  '
$ws.Range("D18").Value = '3.299.20'
$ws.Range("E18").Value = '  -6.40%  '
$ws.Range("D19").Value = '''11.73'
$ws.Range("E19").Value = '  -3.18%  '
$ws.Range("D20").Value = '63.478.14'
$ws.Range("E20").Value = '  -3.70%  '
$ws.Range("D21").Value = '''0.966'
$ws.Range("E21").Value = '  -4.08%  '
$ws.Range("D22").Value = '''427.47'
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").Value = '''4.64'
$ws.Range("E23").Value = '  +7.41%  '
$ws.Range("D24").Value = '''4.07'
$ws.Range("E24").Value = '  -3.48%  '
$ws.Range("D25").Value = '''84.00'
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("D26").Value = '''13.16'
$ws.Range("E26").Value = '  +1.63%  '
$ws.Range("D27").Value = '''10.57'
$ws.Range("E27").Value = '  -3.44%  '
$ws.Range("E28").Value = '  -1.94%  '
$ws.Range("E29").Value = '  -4.22%  '
$ws.Range("D30").Value = '''29.48'
$ws.Range("E30").Value = '  -3.00%  '
$ws.Range("D31").Value = '''6.59'
$ws.Range("E31").Value = '  +2.46%  '
$ws.Range("D32").Value = '''589.48'
$ws.Range("E32").Value = '  -5.75%  '
$ws.Range("D33").Value = '''11.37'
$ws.Range("E33").Value = '  -2.54%  '
$ws.Range("E34").Value = '  -3.76%  '
$ws.Range("D35").Value = '''58.17'
$ws.Range("E35").Value = '  -2.61%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("E37").Value = '  -7.65%  '
$ws.Range("D38").Value = '''35.16'
$ws.Range("E38").Value = '  -5.74%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0742'
$ws.Range("E39").Value = '  -7.60%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''3.39'
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("D41").Value = '''0.363'
$ws.Range("E41").Value = '  -4.41%  '
$ws.Range("D42").Value = '3.083.60'
$ws.Range("E42").Value = '  -6.05%  '
$ws.Range("D43").Value = '''0.999'
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").Value = '''2.79'
$ws.Range("E44").Value = '  -3.62%  '
$ws.Range("D45").Value = '''3.19'
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").Value = '''0.0404'
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("D47").Value = '''2.42'
$ws.Range("E47").Value = '  -3.29%  '
$ws.Range("E48").Value = '  -2.79%  '
$ws.Range("E49").Value = '  -5.47%  '
$ws.Range("D50").Value = '''132.54'
$ws.Range("E50").Value = '  -4.49%  '
$ws.Range("D51").Value = '''8.08'
$ws.Range("E51").Value = '  -5.22%  '
